$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "in Events: Daten für initiale Diagnose und LastVitalStatus ggf. anpassen, falls nicht konsistent mit anderen verfügbaren Daten"
$ws.Range("B6").Value = "Bei InitialDiagnosis: Diagnose-Details zu nestes Details"

$ws.Range("B7").Select()
